$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete the two old "Estado de Cuenta" detail rows (2507, 2508 periods) -
# keep only the first detail row (now effectively representing period 2508)
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(17).Delete()

# Update summary values
$ws.Range("E11").Value = 68000
$ws.Range("F13").Value = 1
